# Updates cryptos list values per the commit diff (Wed Apr 19 09:11:08 UTC 2023 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.412.92'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.81%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.007.03'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -5.08%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.21'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -4.41%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4953'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4180'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -5.97%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.48'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08791'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -6.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.115'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -5.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.221.64'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.14'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -8.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.097'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.478'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -6.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.13'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -6.71%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001102'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -5.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06629'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.50'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -9.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.009'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.967'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.486.47'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.71%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.81'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -7.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.293'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.373.42'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.659'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.26%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.47'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.59'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -6.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.344'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -7.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '127.08'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.049'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -8.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09940'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.83%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -12.02%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.815'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.78%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.791'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.580'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -10.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02453'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06374'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.285'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.82'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -7.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6484'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2064'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -7.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.009'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6310'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -7.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.201'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -6.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '13.38'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -8.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.254'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.561'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07005'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.51%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.143'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.42%  '
